$d = $word.ActiveDocument

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml) | Out-Null
}

# Paragraph 1: "Maxim (Frontend)" heading - drop the en-US language formatting and
# wrap "Frontend" with spell-check proofErr markers, splitting it into 3 runs.
Set-ParagraphXml 1 '<w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Maxim (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frontend</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p></w:body>'

# Paragraph 2: "Animaties ..." body text - merge the three runs into one.
Set-ParagraphXml 2 '<w:body><w:p><w:r><w:t xml:space="preserve">Animaties gemaakt voor het menu en navigatie voor leaderboards is in orde. </w:t></w:r></w:p></w:body>'

# Paragraph 3: "Robbe (Backend)" heading - drop en-US language + spell-check markers,
# merge into a single run.
Set-ParagraphXml 3 '<w:body><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Robbe (Backend)</w:t></w:r></w:p></w:body>'

# Paragraph 4: "Sinds gisteren ..." body text - merge the three runs into one.
Set-ParagraphXml 4 '<w:body><w:p><w:r><w:t xml:space="preserve">Sinds gisteren bugs uit de backend aan het halen. Hij is hier nog mee bezig aangezien er wat meer foutjes tevoorschijn zijn gekomen. </w:t></w:r></w:p></w:body>'

# Paragraph 5: "Herber (Hardware)" heading - unchanged.

# Paragraph 6: "1 led ..." body text - drop the superscript "de" run, merge into one run.
Set-ParagraphXml 6 '<w:body><w:p><w:r><w:t xml:space="preserve">1 led verbonden met de knop en dit werkt perfect. Nu bezig met een 2de laag op de mat te leggen voor de buttons. Dit zou tegen vanavond moeten klaar zijn. </w:t></w:r></w:p></w:body>'

# Paragraph 7: "Jakob (Documentatie)" heading - unchanged.

# Paragraph 8: "Designdocument ..." body text - merge the three runs into one.
Set-ParagraphXml 8 '<w:body><w:p><w:r><w:t>Designdocument is bijna af dit zal vandaag af zijn. Begint dan met ideeën op te zoeken voor de handleidingen.</w:t></w:r></w:p></w:body>'

# Remove the trailing empty paragraph (paragraph 9), which sat right before the sectPr.
$last = $d.Paragraphs($d.Paragraphs.Count)
$r = $d.Range($last.Range.Start - 1, $last.Range.End)
$r.Delete() | Out-Null
